$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Gujarat"
$ws.Range("C18").Value = "Punjab"
$ws.Range("D18").Value = "Punjab"
$ws.Range("E18").Value = 20
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 151

$ws.Range("H18").Select()
